$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data record needs to be inserted as the second-to-last row.
# Insert a blank row at position 54; this pushes the former row 54
# (the last data row) down to row 55, keeping all of its values/styles intact.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new record's data.
$ws.Cells.Item(54, 1).Value = 3
$ws.Cells.Item(54, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(54, 3).Value = "Coquimbo"
$ws.Cells.Item(54, 4).Value = 44628
$ws.Cells.Item(54, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(54, 5).Value = 5
$ws.Cells.Item(54, 6).Value = 100112022
$ws.Cells.Item(54, 7).Value = "Arveja Verde"
$ws.Cells.Item(54, 8).Value = "Perfection"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 73
$ws.Cells.Item(54, 11).Value = 23000
$ws.Cells.Item(54, 12).Value = 24000
$ws.Cells.Item(54, 13).Value = 23521
$ws.Cells.Item(54, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(54, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(54, 16).Value = 941
$ws.Cells.Item(54, 17).Value = 25
$ws.Cells.Item(54, 18).Value = "Hortaliza"
